# Update cryptocurrency symbol list data (Price, Volume(1h), Data, Hora columns)
# for rows 2-51, reflecting the refreshed values from the GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D='276.85'; E='6.13%'; F='9-1-2023'; G='1'},
    @{Row=3; D='27.20'; E='-0.42%'; F='9-1-2023'; G='1'},
    @{Row=4; D='4.769'; E='1.25%'; F='9-1-2023'; G='1'},
    @{Row=5; D='0.06207'; E='-0.85%'; F='9-1-2023'; G='1'},
    @{Row=6; D='6.801'; E='1.26%'; F='9-1-2023'; G='1'},
    @{Row=7; D='0.8653'; E='1.67%'; F='9-1-2023'; G='1'},
    @{Row=8; D='0.9243'; E='1.59%'; F='9-1-2023'; G='1'},
    @{Row=9; D='0.1450'; E='3.94%'; F='9-1-2023'; G='1'},
    @{Row=10; D='0.05184'; E='8.77%'; F='9-1-2023'; G='1'},
    @{Row=11; D='0.07275'; E='2.49%'; F='9-1-2023'; G='1'},
    @{Row=12; D='0.03115'; E='-0.46%'; F='9-1-2023'; G='1'},
    @{Row=13; D='0.09043'; E='-0.10%'; F='9-1-2023'; G='1'},
    @{Row=14; D='0.001544'; E='0.45%'; F='9-1-2023'; G='1'},
    @{Row=15; D='0.0006157'; E='0.03%'; F='9-1-2023'; G='1'},
    @{Row=16; D='0.006002'; E='-2.50%'; F='9-1-2023'; G='1'},
    @{Row=17; D='3.476'; E='0.29%'; F='9-1-2023'; G='1'},
    @{Row=18; D='3.229'; E='1.88%'; F='9-1-2023'; G='1'},
    @{Row=19; D='2.275'; E='5.01%'; F='9-1-2023'; G='1'},
    @{Row=20; D='0.3086'; E='-0.65%'; F='9-1-2023'; G='1'},
    @{Row=21; D='0.1306'; E='0.35%'; F='9-1-2023'; G='1'},
    @{Row=22; D='3.832'; E='-6.16%'; F='9-1-2023'; G='1'},
    @{Row=23; D='0.04249'; E='0.61%'; F='9-1-2023'; G='1'},
    @{Row=24; D='0.001173'; E='-3.87%'; F='9-1-2023'; G='1'},
    @{Row=25; D='0.004216'; E='2.89%'; F='9-1-2023'; G='1'},
    @{Row=26; D='0.0001197'; E='-0.43%'; F='9-1-2023'; G='1'},
    @{Row=27; D='0.0001930'; E='19.40%'; F='9-1-2023'; G='1'},
    @{Row=28; F='9-1-2023'; G='1'},
    @{Row=29; F='9-1-2023'; G='1'},
    @{Row=30; F='9-1-2023'; G='1'},
    @{Row=31; F='9-1-2023'; G='1'},
    @{Row=32; F='9-1-2023'; G='1'},
    @{Row=33; F='9-1-2023'; G='1'},
    @{Row=34; F='9-1-2023'; G='1'},
    @{Row=35; F='9-1-2023'; G='1'},
    @{Row=36; F='9-1-2023'; G='1'},
    @{Row=37; F='9-1-2023'; G='1'},
    @{Row=38; F='9-1-2023'; G='1'},
    @{Row=39; F='9-1-2023'; G='1'},
    @{Row=40; D='0.04015'; E='3.53%'; F='9-1-2023'; G='1'},
    @{Row=41; D='0.006223'; E='52.38%'; F='9-1-2023'; G='1'},
    @{Row=42; D='0.1135'; E='2.16%'; F='9-1-2023'; G='1'},
    @{Row=43; D='0.002116'; E='-3.17%'; F='9-1-2023'; G='1'},
    @{Row=44; D='0.01184'; E='-12.07%'; F='9-1-2023'; G='1'},
    @{Row=45; D='0.00005116'; E='-0.68%'; F='9-1-2023'; G='1'},
    @{Row=46; D='0.00000000748'; E='-0.47%'; F='9-1-2023'; G='1'},
    @{Row=47; D='0.8953'; E='2,732.58%'; F='9-1-2023'; G='1'},
    @{Row=48; D='0.02467'; E='-29.58%'; F='9-1-2023'; G='1'},
    @{Row=49; D='0.00002093'; E='-0.47%'; F='9-1-2023'; G='1'},
    @{Row=50; D='0.0001993'; E='-0.47%'; F='9-1-2023'; G='1'},
    @{Row=51; F='9-1-2023'; G='1'}
)

$columns = @('D', 'E', 'F', 'G')

foreach ($item in $updates) {
    $row = $item.Row
    foreach ($col in $columns) {
        if ($item.ContainsKey($col)) {
            $addr = "$col$row"
            # Force text storage so values are written back as plain strings
            # (matching the original inline-string cell type) instead of being
            # auto-converted to numbers/percentages/dates by Excel.
            $ws.Range($addr).NumberFormat = "@"
            $ws.Range($addr).Value = $item[$col]
            $ws.Range($addr).Style = "Normal"
        }
    }
}

Write-Output "Updated $($updates.Count) rows"
